$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Status" text everywhere it appears: "In Translation" ->
#    "Handed back: in sync with en-US" (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Overview sheet column widths (Status columns grew to fit new text)
# ---------------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# 3. zh-cn sheet: fill in handback report columns (I, J) for both rows,
#    update column widths, and the handback datetime (K)
# ---------------------------------------------------------------------------
$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1").ColumnWidth = 40
$wsZhCn.Range("J1").ColumnWidth = 40

$wsZhCn.Range("I2").Value = "bebe5637-b11f-4808-bd3d-c079eb9e3510.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a12ed0b37bfc539971d0b062f12e29370c87cfa/e2e/bebe5637-b11f-4808-bd3d-c079eb9e3510.md", "", "", "bebe5637-b11f-4808-bd3d-c079eb9e3510.md") | Out-Null
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("J2").Value = "bebe5637-b11f-4808-bd3d-c079eb9e3510.352d601e84afdd5d290fc35ac3ba39cdc7f30f44.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-08 04:35:33"

$wsZhCn.Range("I3").Value = "e01ec07d-238d-4974-b527-b103d0845eab.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a12ed0b37bfc539971d0b062f12e29370c87cfa/e2e/e01ec07d-238d-4974-b527-b103d0845eab.md", "", "", "e01ec07d-238d-4974-b527-b103d0845eab.md") | Out-Null
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("J3").Value = "e01ec07d-238d-4974-b527-b103d0845eab.51f9856749517964d939675fdb759b6e0c175ea3.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-08 04:35:33"

# ---------------------------------------------------------------------------
# 4. de-de sheet: fill in handback report columns (I, J) for both rows,
#    update column widths, and the handback datetime (K)
# ---------------------------------------------------------------------------
$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1").ColumnWidth = 40
$wsDeDe.Range("J1").ColumnWidth = 40

$wsDeDe.Range("I2").Value = "bebe5637-b11f-4808-bd3d-c079eb9e3510.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a12ed0b37bfc539971d0b062f12e29370c87cfa/e2e/bebe5637-b11f-4808-bd3d-c079eb9e3510.md", "", "", "bebe5637-b11f-4808-bd3d-c079eb9e3510.md") | Out-Null
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("J2").Value = "bebe5637-b11f-4808-bd3d-c079eb9e3510.352d601e84afdd5d290fc35ac3ba39cdc7f30f44.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-08 04:35:41"

$wsDeDe.Range("I3").Value = "e01ec07d-238d-4974-b527-b103d0845eab.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a12ed0b37bfc539971d0b062f12e29370c87cfa/e2e/e01ec07d-238d-4974-b527-b103d0845eab.md", "", "", "e01ec07d-238d-4974-b527-b103d0845eab.md") | Out-Null
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("J3").Value = "e01ec07d-238d-4974-b527-b103d0845eab.51f9856749517964d939675fdb759b6e0c175ea3.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-08 04:35:41"

Write-Host "Done applying handback report updates"
